# Rename every "N/A" issue label in column C to "Not Available"
# (applies across the whole "issue" column of the fix-time table),
# and move the active selection to I25 (was I19).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value() -eq "N/A") {
        $cell.Value = "Not Available"
    }
}

[void]$ws.Range("I25").Select()
